$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B4").Value = "levels_CR"
$ws.Range("B5").Value = "levels_fraction"
$ws.Range("B6").Value = "fraction_names"
$ws.Range("B9").Value = "colmns_sorted"

$ws.Range("B3:B6").Interior.Pattern = -4142
$ws.Range("B9").Interior.Pattern = -4142

$ws.Range("D15").Select()
